# Weekly fruit/vegetable price update.
# Inserts a new week's row of data at row 28 (pushing the existing
# rows 28-30 down to 29-31) and populates the new row with the
# latest "Achicoria" price entry for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28, shifting rows 28:30 -> 29:31
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with this week's data
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value = 44694
$ws.Cells.Item(28, 5).Value = 9
$ws.Cells.Item(28, 6).Value = 100112010
$ws.Cells.Item(28, 7).Value = "Achicoria"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 35
$ws.Cells.Item(28, 11).Value = 12000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 12000
$ws.Cells.Item(28, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 667
$ws.Cells.Item(28, 17).Value = 18
$ws.Cells.Item(28, 18).Value = "Hortaliza"
